$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.04382726404007542
$ws.Range("D2").Value = 0.0610941525082076
$ws.Range("E2").Value = 0.07415118220752248
$ws.Range("F2").Value = 1.71459237914236
$ws.Range("G2").Value = 0.0024834632783052
$ws.Range("I2").Value = 0.6485265063505103
$ws.Range("K2").Value = 1.347527313325543
$ws.Range("M2").Value = 0.4200880741095929
$ws.Range("N2").Value = 1.931502099246231

$ws.Range("B3").Value = 0.03850589579796093
$ws.Range("D3").Value = 0.06166449192951617
$ws.Range("E3").Value = 0.06875439041070663
$ws.Range("F3").Value = 1.675056690142142
$ws.Range("G3").Value = 0.002488695114563778
$ws.Range("I3").Value = 0.6499974744972867
$ws.Range("K3").Value = 1.221639574532105
$ws.Range("M3").Value = 0.3832854564961394
$ws.Range("N3").Value = 1.943920516330941

$ws.Range("B4").Value = 0.03523861601355094
$ws.Range("D4").Value = 0.0620359071110439
$ws.Range("E4").Value = 0.06549120452573831
$ws.Range("F4").Value = 1.651884489558114
$ws.Range("G4").Value = 0.002492074519392035
$ws.Range("I4").Value = 0.6511632311583888
$ws.Range("K4").Value = 1.145038627748562
$ws.Range("M4").Value = 0.3609282249870205
$ws.Range("N4").Value = 1.952176144748549

$ws.Range("B5").Value = 0.03390732221785697
$ws.Range("D5").Value = 0.06219257651186716
$ws.Range("E5").Value = 0.06417390214585694
$ws.Range("F5").Value = 1.642717283082021
$ws.Range("G5").Value = 0.002493493805350792
$ws.Range("I5").Value = 0.6517041993856196
$ws.Range("K5").Value = 1.113995168526031
$ws.Range("M5").Value = 0.3518767987256695
$ws.Range("N5").Value = 1.95569856305022

$ws.Range("B6").Value = 0.03368627485681941
$ws.Range("D6").Value = 0.0622189116950338
$ws.Range("E6").Value = 0.06395591206518603
$ws.Range("F6").Value = 1.641211675710395
$ws.Range("G6").Value = 0.002493732026938456
$ws.Range("I6").Value = 0.6517980044146761
$ws.Range("K6").Value = 1.108850741217339
$ws.Range("M6").Value = 0.3503773718711756
$ws.Range("N6").Value = 1.956293000024232

$ws.Range("B7").Value = 0.03522066093748322
$ws.Range("D7").Value = 0.06203799852105973
$ws.Range("E7").Value = 0.06547338870975494
$ws.Range("F7").Value = 1.651759743516962
$ws.Range("G7").Value = 0.002492093489482843
$ws.Range("I7").Value = 0.6511702601226474
$ws.Range("K7").Value = 1.144619271622219
$ws.Range("M7").Value = 0.3608059153230201
$ws.Range("N7").Value = 1.95222300939421

$ws.Range("B8").Value = 0.0419925217535706
$ws.Range("D8").Value = 0.06128638686993959
$ws.Range("E8").Value = 0.07227976086112875
$ws.Range("F8").Value = 1.700730543973336
$ws.Range("G8").Value = 0.002485232636656392
$ws.Range("I8").Value = 0.6489791286758013
$ws.Range("K8").Value = 1.303975421756149
$ws.Range("M8").Value = 0.4073481881950514
$ws.Range("N8").Value = 1.935652829675753

$ws.Range("B9").Value = 0.05526740356897619
$ws.Range("D9").Value = 0.05998188823687656
$ws.Range("E9").Value = 0.08603804604451426
$ws.Range("F9").Value = 1.805595996569849
$ws.Range("G9").Value = 0.002473097095746338
$ws.Range("I9").Value = 0.6467713504177084
$ws.Range("K9").Value = 1.622119907275078
$ws.Range("M9").Value = 0.5005677235905637
$ws.Range("N9").Value = 1.908181124179436

$ws.Range("B10").Value = 0.06501156416595677
$ws.Range("D10").Value = 0.05912813054542099
$ws.Range("E10").Value = 0.09641257847457751
$ws.Range("F10").Value = 1.888152023048804
$ws.Range("G10").Value = 0.002464975429199833
$ws.Range("I10").Value = 0.6464312359527895
$ws.Range("K10").Value = 1.859517296895262
$ws.Range("M10").Value = 0.5703192176193994
$ws.Range("N10").Value = 1.89108481355305

$ws.Range("B11").Value = 0.06944129654746689
$ws.Range("D11").Value = 0.05876276722043983
$ws.Range("E11").Value = 0.1011935519063982
$ws.Range("F11").Value = 1.92693385348295
$ws.Range("G11").Value = 0.002461451114710695
$ws.Range("I11").Value = 0.6465567562988142
$ws.Range("K11").Value = 1.968357763462222
$ws.Range("M11").Value = 0.6023418855721872
$ws.Range("N11").Value = 1.883983210382453

$ws.Range("B12").Value = 0.07111817771490792
$ws.Range("D12").Value = 0.05862775055177138
$ws.Range("E12").Value = 0.1030130990231868
$ws.Range("F12").Value = 1.94179804131636
$ws.Range("G12").Value = 0.002460140877374609
$ws.Range("I12").Value = 0.6466447325975153
$ws.Range("K12").Value = 2.009698104824338
$ws.Range("M12").Value = 0.614511298886498
$ws.Range("N12").Value = 1.881391651815562

$ws.Range("B13").Value = 0.07075705858738957
$ws.Range("D13").Value = 0.05865667985820266
$ws.Range("E13").Value = 0.1026208184864004
$ws.Range("F13").Value = 1.938588808866257
$ws.Range("G13").Value = 0.002460421979989367
$ws.Range("I13").Value = 0.6466239841966939
$ws.Range("K13").Value = 2.000789123069126
$ws.Range("M13").Value = 0.6118884626194045
$ws.Range("N13").Value = 1.881945438377983

$ws.Range("B14").Value = 0.06957926664075842
$ws.Range("D14").Value = 0.05875159224123294
$ws.Range("E14").Value = 0.10134306341957
$ws.Range("F14").Value = 1.928153153448164
$ws.Range("G14").Value = 0.002461342833505244
$ws.Range("I14").Value = 0.6465631827498655
$ws.Range("K14").Value = 1.971756339339663
$ws.Range("M14").Value = 0.6033421999983091
$ws.Range("N14").Value = 1.883768040982375

$ws.Range("B15").Value = 0.06885775768263613
$ws.Range("D15").Value = 0.05881016445728804
$ws.Range("E15").Value = 0.1005615935027109
$ws.Range("F15").Value = 1.921784298207115
$ws.Range("G15").Value = 0.002461910049044926
$ws.Range("I15").Value = 0.6465312115460762
$ws.Range("K15").Value = 1.953989265820496
$ws.Range("M15").Value = 0.5981130158115064
$ws.Range("N15").Value = 1.884897171911263

$ws.Range("B16").Value = 0.06472199966461289
$ws.Range("D16").Value = 0.05915247348062991
$ws.Range("E16").Value = 0.09610139003599016
$ws.Range("F16").Value = 1.885642389604953
$ws.Range("G16").Value = 0.002465209165617638
$ws.Range("I16").Value = 0.6464286845096225
$ws.Range("K16").Value = 1.852421618875098
$ws.Range("M16").Value = 0.56823243619084
$ws.Range("N16").Value = 1.891562558522182

$ws.Range("B17").Value = 0.06218399423893572
$ws.Range("D17").Value = 0.05936838410789491
$ws.Range("E17").Value = 0.09338113714078844
$ws.Range("F17").Value = 1.863786063052004
$ws.Range("G17").Value = 0.002467276577505571
$ws.Range("I17").Value = 0.6464376671610665
$ws.Range("K17").Value = 1.79033225830301
$ws.Range("M17").Value = 0.5499773144032503
$ws.Range("N17").Value = 1.895824956873284

$ws.Range("B18").Value = 0.06072393270696352
$ws.Range("D18").Value = 0.05949473473069489
$ws.Range("E18").Value = 0.09182229777164963
$ws.Range("F18").Value = 1.851330224288233
$ws.Range("G18").Value = 0.002468481733459305
$ws.Range("I18").Value = 0.6464692024139325
$ws.Range("K18").Value = 1.754699660919982
$ws.Range("M18").Value = 0.5395049040337625
$ws.Range("N18").Value = 1.898340134538174

$ws.Range("B19").Value = 0.06022953922990837
$ws.Range("D19").Value = 0.05953788577157582
$ws.Range("E19").Value = 0.09129548651124253
$ws.Range("F19").Value = 1.847132642340682
$ws.Range("G19").Value = 0.002468892536633022
$ws.Range("I19").Value = 0.6464844041234485
$ws.Range("K19").Value = 1.742648650023625
$ws.Range("M19").Value = 0.5359638081822453
$ws.Range("N19").Value = 1.899202631209747

$ws.Range("B20").Value = 0.06245419799844854
$ws.Range("D20").Value = 0.05934517582393184
$ws.Range("E20").Value = 0.09367011291846694
$ws.Range("F20").Value = 1.86610075478302
$ws.Range("G20").Value = 0.002467054839445174
$ws.Range("I20").Value = 0.6464339808330379
$ws.Range("K20").Value = 1.796933521006338
$ws.Range("M20").Value = 0.5519177542300753
$ws.Range("N20").Value = 1.895364635490964

$ws.Range("B21").Value = 0.06992522906173804
$ws.Range("D21").Value = 0.05872362333446546
$ws.Range("E21").Value = 0.1017181222100589
$ws.Range("F21").Value = 1.931213502296259
$ws.Range("G21").Value = 0.002461071697155469
$ws.Range("I21").Value = 0.6465799427813934
$ws.Range("K21").Value = 1.980280560026074
$ws.Range("M21").Value = 0.6058512654097683
$ws.Range("N21").Value = 1.883230043686979

$ws.Range("B22").Value = 0.07480463744850852
$ws.Range("D22").Value = 0.05833687606621751
$ws.Range("E22").Value = 0.1070310979467521
$ws.Range("F22").Value = 1.97480893835106
$ws.Range("G22").Value = 0.002457303199637511
$ws.Range("I22").Value = 0.6469111479733201
$ws.Range("K22").Value = 2.100837415900742
$ws.Range("M22").Value = 0.641351807058939
$ws.Range("N22").Value = 1.875868894302016

$ws.Range("B23").Value = 0.07220075731514441
$ws.Range("D23").Value = 0.05854149882147652
$ws.Range("E23").Value = 0.1041905186734482
$ws.Range("F23").Value = 1.951445374565395
$ws.Range("G23").Value = 0.002459301586290985
$ws.Range("I23").Value = 0.6467127519818163
$ws.Range("K23").Value = 2.036426249802446
$ws.Range("M23").Value = 0.6223810911069307
$ws.Range("N23").Value = 1.879745395988508

$ws.Range("B24").Value = 0.06233204176022866
$ws.Range("D24").Value = 0.05935566137424431
$ws.Range("E24").Value = 0.09353945119514862
$ws.Range("F24").Value = 1.865053941249556
$ws.Range("G24").Value = 0.002467155035467631
$ws.Range("I24").Value = 0.6464355652931602
$ws.Range("K24").Value = 1.79394889287687
$ws.Range("M24").Value = 0.5510404107531883
$ws.Range("N24").Value = 1.895572545400071

$ws.Range("B25").Value = 0.05167732682984649
$ws.Range("D25").Value = 0.06031651294392759
$ws.Range("E25").Value = 0.08227044632079838
$ws.Range("F25").Value = 1.776268464166975
$ws.Range("G25").Value = 0.002476239905276001
$ws.Range("I25").Value = 0.6471441054161104
$ws.Range("K25").Value = 1.535426873287747
$ws.Range("M25").Value = 0.4751329715155634
$ws.Range("N25").Value = 1.915072674647078
